$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Capacity")

# Update F35 value (21000 -> 25000)
$ws.Range("F35").Value = 25000

# New rows 47-50
$ws.Range("A47").Value = "EE00"
$ws.Range("C47").Value = "Onshore Wind"
$ws.Range("D47").Value = "Distributed Energy"
$ws.Range("E47").Value = 2040
$ws.Range("F47").Value = 2500

$ws.Range("A48").Value = "LV00"
$ws.Range("C48").Value = "Onshore Wind"
$ws.Range("D48").Value = "Distributed Energy"
$ws.Range("E48").Value = 2040
$ws.Range("F48").Value = 2500

$ws.Range("A49").Value = "LT00"
$ws.Range("C49").Value = "Onshore Wind"
$ws.Range("D49").Value = "Distributed Energy"
$ws.Range("E49").Value = 2040
$ws.Range("F49").Value = 2500

$ws.Range("A50").Value = "UK00"
$ws.Range("C50").Value = "Onshore Wind"
$ws.Range("D50").Value = "Distributed Energy"
$ws.Range("E50").Value = 2040
$ws.Range("F50").Value = 25000

# Selection change to C47
$ws.Range("C47").Select()
